# Apply the "purpose" column correction: S.GISH -> fullRNASEQ for rows 2-21
# (the preparer's name had mistakenly been duplicated into the purpose
# column; this restores the real purpose value, adding a new shared
# string "fullRNASEQ" in the process)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

for ($r = 2; $r -le 21; $r++) {
    $ws.Cells.Item($r, 5).Value = "fullRNASEQ"
}

# Restore the scroll position / selection state recorded in the sheet view
$win = $excel.ActiveWindow
$win.ScrollRow = 13
$win.ScrollColumn = 1
$ws.Range("D22:F24").Select()

# Enable iterative calculation with a max change of 1E-4
$wb.Application.Iteration = $true
$wb.Application.MaxChange = 0.0001
